$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- "carga de usuarios via csv o excel" -----------------------------------
# The "telefono" column (C) is no longer needed/used by the CSV/Excel user
# loader, so drop it entirely. "email" (was D) and "password" (was E) shift
# left into C and D.

# The existing mailto: hyperlink lives on the "email" cell of the 3rd
# student row (currently D4). Stash its current (custom, non-default)
# formatting in a scratch cell so we can restore it byte-for-byte after the
# column shuffle - otherwise Excel would stamp the cell with its own
# built-in "Hyperlink" style.
$ws.Range("D4").Copy()
$ws.Range("ZZ1").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

# Hyperlink anchors are not re-targeted automatically when the sheet shifts
# underneath them, so clear the old one before moving data around.
$ws.Hyperlinks.Delete()

# Remove the whole "telefono" column - matricula/email/password shift left.
$ws.Columns("C").Delete()

# Re-create the hyperlink on its new home cell, C4 (still the 3rd student's
# email cell, just one column to the left now).
$ws.Hyperlinks.Add($ws.Range("C4"), "mailto:estudiante3@pucmm.edu.do", [Type]::Missing, [Type]::Missing, "estudiante3@pucmm.edu.do")

# Restore the original cell formatting (the scratch cell shifted from ZZ1 to
# ZY1 along with everything else once column C was deleted) and clean up.
$ws.Range("ZY1").Copy()
$ws.Range("C4").PasteSpecial(-4122)
$excel.CutCopyMode = $false
$ws.Range("ZY1").Clear()

# Adding a hyperlink registers Excel's built-in "Hyperlink" cell style even
# though we don't actually use it on any cell (we restored the original
# formatting above) - drop it again so the style table stays as it was.
$wb.Styles("Hyperlink").Delete()

# Match the cursor position left behind by the edit.
$ws.Range("E6").Select()
